# Append two new paragraphs at the end of the document:
#   1) "Version management" + ":"   (as two separate runs)
#   2) the long descriptive paragraph about version management
#
# The existing trailing empty paragraph is left empty/untouched in content.

$d = $word.ActiveDocument

# --- Paragraph 1: "Version management" + ":" --------------------------
# Insert a new empty paragraph after the very end of the document.
$endRange = $d.Content
$endRange.Collapse(0)
$endRange.InsertParagraphAfter() | Out-Null

# Type "Version management" into that new paragraph, then split it into
# its own paragraph so the trailing ":" can be typed as a distinct run.
$paraVM = $d.Paragraphs.Last
$rangeVM = $paraVM.Range
$rangeVM.InsertAfter("Version management") | Out-Null
$rangeVM.Collapse(0)
$rangeVM.InsertParagraphAfter() | Out-Null

# Type the ":" into the newly split-off paragraph.
$paraColon = $d.Paragraphs.Last
$rangeColon = $paraColon.Range
$rangeColon.InsertAfter(":") | Out-Null

# Re-join the two paragraphs by removing the paragraph mark between them,
# leaving "Version management" and ":" as two separate runs inside one
# paragraph.
$joinPara = $d.Paragraphs.Item($d.Paragraphs.Count - 1)
$joinRange = $joinPara.Range
$joinRange.Collapse(0)
$joinRange.MoveEnd(1, 1) | Out-Null
$joinRange.Delete() | Out-Null

# --- Paragraph 2: the long descriptive paragraph -----------------------
$mergedPara = $d.Paragraphs.Last
$afterRange = $mergedPara.Range
$afterRange.Collapse(0)
$afterRange.InsertParagraphAfter() | Out-Null

$paraBody = $d.Paragraphs.Last
$rangeBody = $paraBody.Range
$rangeBody.InsertAfter("As a system goes through the development, system testing and release phases discussed in the previous section, each component has many versions.  The components can be code files, configuration files, data file, documentation, media items, or any digital item that is required to build an entire system. Version management is about managing these many components.  Since many people may be working on the components at the same time, version management has to manage these components so that they do not interfere with each other.") | Out-Null
